$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.620500000000006
$ws.Range("B9").Value = 5.277700000000005
$ws.Range("D9").Value = -8.38559999999999
$ws.Range("B18").Value = 6.125899999999998
$ws.Range("B20").Value = 9.124600000000003
$ws.Range("D23").Value = -7.945599999999998
$ws.Range("D24").Value = -7.288599999999998
$ws.Range("D26").Value = -7.592800000000004
$ws.Range("B27").Value = 6.323800000000003
$ws.Range("D34").Value = -7.831700000000001
$ws.Range("D35").Value = -7.460300000000003
$ws.Range("D48").Value = -7.381499999999997
$ws.Range("D52").Value = -7.776600000000002
$ws.Range("D66").Value = -7.1481
$ws.Range("D67").Value = -6.873599999999998
$ws.Range("B69").Value = 5.605999999999989
$ws.Range("B76").Value = 5.062699999999998
$ws.Range("D80").Value = -8.026700000000003
$ws.Range("B82").Value = 5.481
$ws.Range("D99").Value = -8.174400000000002
